# Fruta / hortaliza, semanal
# Inserts a new weekly price record for "Feria Lagunitas de Puerto Montt - Membrillo"
# as a new row 75, pushing the previously-existing rows 75-90 down to 76-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 75 (shifts 75..90 -> 76..91)
$ws.Rows("75").Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44722
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100104
$ws.Range("H75").Value = "Frutos de pepita"
$ws.Range("I75").Value = 100104003
$ws.Range("J75").Value = "Membrillo"
$ws.Range("K75").Value = "Champion"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 300
$ws.Range("N75").Value = 13000
$ws.Range("O75").Value = 14000
$ws.Range("P75").Value = 13500
$ws.Range("Q75").Value = "$/caja 18 kilos granel"
$ws.Range("R75").Value = "Región de O'Higgins"
$ws.Range("S75").Value = 750
$ws.Range("T75").Value = 18
